$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1366
$ws.Range("K5").Value = 718
$ws.Range("K6").Value = 648
